$wb = $excel.ActiveWorkbook

# --- Sheet: "Totales Plantel 2P" ---
$ws = $wb.Worksheets.Item("Totales Plantel 2P")
$ws.Range("F3").Value = 17
$ws.Range("H3").Value = 14
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 41
$ws.Range("K3").Value = 82
$ws.Range("L3").Value = 119
$ws.Range("M3").Value = 59.2

$ws.Range("H4").Value = 14
$ws.Range("I4").Value = 26

$ws.Range("H7").Value = 11
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 12
$ws.Range("K7").Value = 53
$ws.Range("L7").Value = 72
$ws.Range("M7").Value = 57.6

# --- Sheet: "Totales Plantel Final" ---
$ws = $wb.Worksheets.Item("Totales Plantel Final")
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 19
$ws.Range("J3").Value = 33
$ws.Range("K3").Value = 99
$ws.Range("L3").Value = 102
$ws.Range("M3").Value = 50.75

$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 115
$ws.Range("L4").Value = 74
$ws.Range("M4").Value = 39.15

# --- Sheet: "Reprobados por Grupo" ---
$ws = $wb.Worksheets.Item("Reprobados por Grupo")
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 2
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 13
$ws.Range("K10").Value = 15
$ws.Range("L10").Value = 53.57

$ws.Range("I12").Value = 7
$ws.Range("J12").Value = 12
$ws.Range("K12").Value = 25
$ws.Range("L12").Value = 67.56999999999999

$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 4

$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 27
$ws.Range("K17").Value = 7
$ws.Range("L17").Value = 20.59

# --- Sheet: "Totales Grupos" ---
$ws = $wb.Worksheets.Item("Totales Grupos")
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 42.86
$ws.Range("G10").Value = 13
$ws.Range("H10").Value = 46.43

$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 21.62
$ws.Range("G12").Value = 12
$ws.Range("H12").Value = 32.43

$ws.Range("E13").Value = 21
$ws.Range("F13").Value = 58.33

$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 79.41

$ws.Range("E30").Value = 11
$ws.Range("F30").Value = 37.93
